$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 updates
$ws.Range("G8").Value = 2.75
$ws.Range("H8").Value = 2.88
$ws.Range("I8").Value = 2.8
$ws.Range("T8").Value = 6.5
$ws.Range("U8").Value = 12
$ws.Range("V8").Value = 11
$ws.Range("W8").Value = 29
$ws.Range("X8").Value = 29
$ws.Range("AC8").Value = 67
$ws.Range("AE8").Value = 7
$ws.Range("AF8").Value = 12
$ws.Range("AG8").Value = 12
$ws.Range("AH8").Value = 29
$ws.Range("AI8").Value = 29
$ws.Range("AJ8").Value = 41

# Row 9 updates
$ws.Range("P9").Value = 1.8
$ws.Range("Q9").Value = 2
